$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.799.73"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3
$ws.Range("D3").Value = "2.775.19"
$ws.Range("E3").Value = "  -1.55%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "'355.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "

# Row 6
$ws.Range("D6").Value = "'109.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.75%  "

# Row 7
$ws.Range("E7").Value = "  +1.72%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.97%  "

# Row 10
$ws.Range("E10").Value = "  -3.86%  "

# Row 11
$ws.Range("D11").Value = "'0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("D13").Value = "'19.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.91%  "

# Row 14
$ws.Range("D14").Value = "'7.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

# Row 15
$ws.Range("D15").Value = "3.204.64"
$ws.Range("E15").Value = "  -2.10%  "

# Row 16
$ws.Range("D16").Value = "2.788.62"
$ws.Range("E16").Value = "  -1.96%  "

# Row 17
$ws.Range("E17").Value = "  +4.17%  "

# Row 18
$ws.Range("D18").Value = "51.717.50"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("D19").Value = "'7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.86%  "

# Row 20
$ws.Range("E20").Value = "  +0.59%  "

# Row 21
$ws.Range("E21").Value = "  -3.55%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23
$ws.Range("D23").Value = "'274.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.67%  "

# Row 24
$ws.Range("D24").Value = "'69.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "

# Row 25
$ws.Range("E25").Value = "  -1.81%  "

# Row 26
$ws.Range("D26").Value = "'26.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -1.37%  "

# Row 29
$ws.Range("E29").Value = "  +3.75%  "

# Row 30
$ws.Range("E30").Value = "  -1.34%  "

# Row 31
$ws.Range("D31").Value = "'51.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.06%  "

# Row 32
$ws.Range("D32").Value = "'0.0463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.04%  "

# Row 33
$ws.Range("D33").Value = "'33.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "

# Row 34
$ws.Range("E34").Value = "  -1.96%  "

# Row 35
$ws.Range("D35").Value = "'0.0844"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.16%  "

# Row 36
$ws.Range("D36").Value = "'5.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.92%  "

# Row 37
$ws.Range("E37").Value = "  -0.20%  "

# Row 38
$ws.Range("E38").Value = "  +0.19%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'18.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "

# Row 41
$ws.Range("D41").Value = "'0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "

# Row 42
$ws.Range("D42").Value = "'2.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.92%  "

# Row 43
$ws.Range("E43").Value = "  -2.66%  "

# Row 44
$ws.Range("D44").Value = "'121.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.59%  "

# Row 45
$ws.Range("D45").Value = "'21.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.97%  "

# Row 46
$ws.Range("D46").Value = "2.062.01"
$ws.Range("E46").Value = "  -0.95%  "

# Row 47
$ws.Range("D47").Value = "'3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.38%  "

# Row 48
$ws.Range("D48").Value = "'2.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.84%  "

# Row 49
$ws.Range("E49").Value = "  +0.18%  "

# Row 50
$ws.Range("D50").Value = "'0.927"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "

# Row 51
$ws.Range("E51").Value = "  +0.28%  "

